$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete "(population census results)" caption text and the
# now-unused helper cells around it / above the title row.
$ws.Range("A2").Clear() | Out-Null
$ws.Range("B1:B3").Clear() | Out-Null

# The row that only held that stray styled cell is now completely empty -
# delete it so the table below moves back up.
$ws.Rows(3).Delete() | Out-Null

# Drop the 1989 and 2002 columns - only one year column (now showing 2014)
# is kept.
$ws.Columns("C:D").Delete() | Out-Null

# The remaining year column used to read 1989; it now represents 2014.
$ws.Range("B4").Value = 2014

# Rename the sheet to match the municipality.
$ws.Name = "ხულო"

# Restore the cursor to A2, matching the saved selection.
$ws.Range("A2").Select() | Out-Null
